$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spells")

# --- Header row ---
$ws.Range("D1").Value = "Desc"

# --- Existing spell rows: fill in new "Desc2" text / tweak values ---
$ws.Range("E2").Value = "TBD"
$ws.Range("E4").Value = "TBD"

$ws.Range("D6").Value = "Deal 5 damage to a single body part"
$ws.Range("E6").Value = "Deal 3 damage to two different body parts"

# Flame Gambit (row 7) reworked
$ws.Range("B7").Value = "F"
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = "For each stack, name a number 1-6. Roll a D6. If it equals one of your numbers, deal damage equal to the number of stacks."
$ws.Range("E7").Value = "Place a marker on Casino for each stack. Your next Casino roll you may modify the roll +/- 1 for each marker."

# --- New spell rows ---
$ws.Range("A8").Value = "Squall"
$ws.Range("B8").Value = "AW"
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "Deal 1 damage"
$ws.Range("E8").Value = "Extinguish 1 Flame"

$ws.Range("A9").Value = "Hurricane"
$ws.Range("B9").Value = "AAWWW"
$ws.Range("D9").Value = "Deal 4 damage"
$ws.Range("E9").Value = "Deal 3 damage. Deal 1 damage for each player who contributes an AW to this spell."

$ws.Range("A10").Value = "Rock Shield"
$ws.Range("B10").Value = "PSS"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "Prevent 2 damage"
$ws.Range("E10").Value = "Deal 2 damage"

# --- Formatting: D/E columns become wide, left-aligned text columns ---
$ws.Columns("D").ColumnWidth = 56.17
$ws.Range("D1:E1").Font.Bold = $true
$ws.Range("D1:E1").HorizontalAlignment = -4131
$ws.Range("D2:E10").HorizontalAlignment = -4131

# --- Selection, matching the saved workbook view ---
[void]$ws.Range("D7").Select()
